$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("stats")

$values = @{
    "D2" = 0.0000930004753172397613525
    "E2" = 0.0293510430492460692997181
    "G2" = 0.001986057497560977935791
    "H2" = 0.003223590552806853814366
    "I2" = 0.0090468446724116802215576
    "J2" = 0.0130320801399648206891912
    "K2" = 0.0005951202474534511566162
    "D3" = 0.0018217079341411590576172
    "E3" = 0.0384883820079267024993896
    "G3" = 0.0021105818450450901553606
    "H3" = 0.0052878647111356258392334
    "I3" = 0.0106542510911822301683527
    "J3" = 0.0181085243821144104003906
    "K3" = 0.0006697354838252067565918
    "D4" = 0.0018931427039206030153023
    "E4" = 0.0408920291811227798461914
    "G4" = 0.0021542995236814022064209
    "H4" = 0.0054485131986439228057861
    "I4" = 0.0103793251328170299530029
    "J4" = 0.0203893352299928699855602
    "K4" = 0.0008129198104143142700195
    "D5" = 0.0002007568255066871914117
    "E5" = 0.0332169309258460998535156
    "G5" = 0.0020764674991369251773332
    "H5" = 0.0037319352850317959353899
    "I5" = 0.0104106073267757892608643
    "J5" = 0.0146356159821152704419989
    "K5" = 0.0007260786369442939758301
    "D6" = 0.0032465481199324131011963
    "E6" = 0.2083004242740570943759337
    "G6" = 0.0055916714482009410858154
    "H6" = 0.0142188328318297897701061
    "I6" = 0.1558146933093667030334473
    "J6" = 0.0256038429215550387973988
    "K6" = 0.0021174303255975250766252
    "D8" = 0.0000930004753172397613525
    "E8" = 0.0293510430492460692997181
    "G8" = 0.001986057497560977935791
    "H8" = 0.003223590552806853814366
    "I8" = 0.0090468446724116802215576
    "J8" = 0.0130320801399648206891912
    "K8" = 0.0005951202474534511566162
    "D9" = 0.0018217079341411590576172
    "E9" = 0.0384883820079267024993896
    "G9" = 0.0021105818450450901553606
    "H9" = 0.0052878647111356258392334
    "I9" = 0.0106542510911822301683527
    "J9" = 0.0181085243821144104003906
    "K9" = 0.0006697354838252067565918
    "D10" = 0.0018931427039206030153023
    "E10" = 0.0408920291811227798461914
    "G10" = 0.0021542995236814022064209
    "H10" = 0.0054485131986439228057861
    "I10" = 0.0103793251328170299530029
    "J10" = 0.0203893352299928699855602
    "K10" = 0.0008129198104143142700195
    "D11" = 0.0002007568255066871914117
    "E11" = 0.0332169309258460998535156
    "G11" = 0.0020764674991369251773332
    "H11" = 0.0037319352850317959353899
    "I11" = 0.0104106073267757892608643
    "J11" = 0.0146356159821152704419989
    "K11" = 0.0007260786369442939758301
    "D12" = 0.0032465481199324131011963
    "E12" = 0.2083004242740570943759337
    "G12" = 0.0055916714482009410858154
    "H12" = 0.0142188328318297897701061
    "I12" = 0.1558146933093667030334473
    "J12" = 0.0256038429215550387973988
    "K12" = 0.0021174303255975250766252
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
